$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 271; this shifts current rows 271..332 down to 272..333
$ws.Rows.Item(271).Insert()

# Populate the newly inserted row 271 with its data
$ws.Cells.Item(271, 1).Value = 3
$ws.Cells.Item(271, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(271, 3).Value = "Coquimbo"
$ws.Cells.Item(271, 4).Value = 44782
$ws.Cells.Item(271, 4).Style = $ws.Cells.Item(272, 4).Style
$ws.Cells.Item(271, 4).NumberFormat = $ws.Cells.Item(272, 4).NumberFormat
$ws.Cells.Item(271, 5).Value = 5
$ws.Cells.Item(271, 6).Value = 100112001
$ws.Cells.Item(271, 7).Value = "Berenjena"
$ws.Cells.Item(271, 8).Value = "Sin especificar"
$ws.Cells.Item(271, 9).Value = "Primera"
$ws.Cells.Item(271, 10).Value = 105
$ws.Cells.Item(271, 11).Value = 8500
$ws.Cells.Item(271, 12).Value = 9000
$ws.Cells.Item(271, 13).Value = 8762
$ws.Cells.Item(271, 14).Value = "`$/caja 60 unidades"
$ws.Cells.Item(271, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(271, 16).Value = 146
$ws.Cells.Item(271, 17).Value = 60
$ws.Cells.Item(271, 18).Value = "Hortaliza"
